$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.606.86"
$ws.Range("E2").Value = "  +0.49%  "

$ws.Range("D3").Value = "2.143.63"
$ws.Range("E3").Value = "  +1.90%  "

$ws.Range("E4").Value = "  +0.33%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.19%  "

$ws.Range("E6").Value = "  +0.24%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5262"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.89%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4566"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.36%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.73"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09181"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.185"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.49%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.61%  "

$ws.Range("D13").Value = "2.141.05"
$ws.Range("E13").Value = "  +1.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.888"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.169"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "102.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.49%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001174"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.83%  "

$ws.Range("E18").Value = "  +0.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06710"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.25%  "

$ws.Range("E21").Value = "  +0.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.356"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.73%  "

$ws.Range("D23").Value = "30.707.09"
$ws.Range("E23").Value = "  +0.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.382"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.00%  "

$ws.Range("D26").Value = "2.364.43"
$ws.Range("E26").Value = "  +0.62%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.95%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.656"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.97%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "164.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.223"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.64%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1082"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.675"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.384"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.020"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.90%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.180"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.59%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02658"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06988"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.32%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2343"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.88%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7031"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.40%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.278"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.55%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.87%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.367"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6483"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.04%  "

$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.758"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.57%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000367"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.29%  "

$ws.Range("E49").Value = "  +0.75%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07307"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.60%  "
